$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A3: was "RO.ACT" -> now "AD.SEC.002.FON.01"
$ws.Range("A3").Value = "AD.SEC.002.FON.01"

# Clear A4 (previously held "AD.SEC.002.FON.01")
$ws.Range("A4").ClearContents()

# Update selection to A8
$ws.Range("A8").Select()
